# Added Data Cleansing Notebook
# Adds a new "Sheet2" after the existing "Sheet1" holding the cleansed
# SharePoint list data (ID, Priority, Company, HQ, Included, Notes_1, Notes_2).

$wb = $excel.ActiveWorkbook

# Insert the new worksheet immediately after Sheet1
$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "Priority"
$ws2.Range("C1").Value = "Company"
$ws2.Range("D1").Value = "HQ"
$ws2.Range("E1").Value = "Included"
$ws2.Range("F1").Value = "Notes_1"
$ws2.Range("G1").Value = "Notes_2"

# Data rows
$data = @(
    @(1,  1, "Company 001", "Ireland",   1, 0.87001607940156156, 0.18366716439731878),
    @(2,  1, "Company 002", "USA",       0, 0.09238234628081321, 0.11574112636337619),
    @(3,  2, "Company 003", "Hong Kong", 1, 0.27007960504501938, 0.19845950723333794),
    @(4,  3, "Company 004", "Ireland",   0, 0.59801876806207455, 0.58131423496351209),
    @(5,  3, "Company 005", "USA",       1, 0.17695636781258961, 0.35726218281172617),
    @(6,  2, "Company 006", "Hong Kong", 1, 0.93497349335220259, 0.86843209899135509),
    @(7,  2, "Company 007", "USA",       0, 0.07789386332132597, 0.31234413831136498),
    @(8,  1, "Company 008", "USA",       1, 0.27333504899498684, 0.51202917172740459),
    @(9,  3, "Company 009", "Ireland",   0, 0.700316980058725,   0.468006308488357),
    @(9,  3, "Company 009", "Ireland",   0, 0.700316980058725,   0.468006308488357),
    @(10, 1, "Company 010", "Hong Kong", 1, 0.2591384043226258,  0.5554633890754721),
    @(10, 1, "Company 010", "Hong Kong", 1, 0.2591384043226258,  0.5554633890754721),
    @(10, 1, "Company 010", "Hong Kong", 1, 0.2591384043226258,  0.5554633890754721)
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# Column widths to (best-) fit the Company / HQ columns
$ws2.Columns.Item(3).ColumnWidth = 13.1
$ws2.Columns.Item(4).ColumnWidth = 23.33

# Make Sheet2 the active sheet/view: zoom 200%, selection on D14
$ws2.Activate()
$excel.ActiveWindow.Zoom = 200
$ws2.Range("D14").Select()
